$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws.Range("B9").Value = "Alvearie Team"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
$ws.Range("A11").Value = "Description"
$ws.Range("B11").Value = "IBM® Health Data Connect indicator type value set"
$ws.Range("A12").Value = "Purpose"
$ws.Range("B12").Value = ""
$ws.Range("A13").Value = "Copyright"
$ws.Range("B13").Value = ""
$ws.Range("A14").Value = "Immutable"
$ws.Range("B14").Value = "BooleanType[null]"

$ws.Rows.Item(15).Delete()
